$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Serial=5, Make=Fitbit, Model=Charge HR, Size=Small, In Possession=Larissa,
#        Position @ Center=Summer Student, Out Date=Aug 29th, 2017
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "Fitbit"
$ws.Range("C4").Value = "Charge HR"
$ws.Range("D4").Value = "Small"
$ws.Range("E4").Value = "Larissa"
$ws.Range("F4").Value = "Summer Student"
$ws.Range("G4").Value = "Aug 29th, 2017"

# Row 5: Serial=1, Make=Fitbit, Model=Charge HR, Size=Small, In Possession=Larissa,
#        Position @ Center=Summer Student, Out Date=Aug 29th, 2017
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Fitbit"
$ws.Range("C5").Value = "Charge HR"
$ws.Range("D5").Value = "Small"
$ws.Range("E5").Value = "Larissa"
$ws.Range("F5").Value = "Summer Student"
$ws.Range("G5").Value = "Aug 29th, 2017"

# Selection moved to E8 in the saved view state.
$ws.Range("E8").Select() | Out-Null
